$wb = $excel.ActiveWorkbook

# --- Step 1: capture the original Sheet1 (has the pre-existing data) and add 5 new sheets after it ---
$old = $wb.Worksheets.Item(1)
$n2 = $wb.Worksheets.Add($null, $old)
$n3 = $wb.Worksheets.Add($null, $n2)
$n4 = $wb.Worksheets.Add($null, $n3)
$n5 = $wb.Worksheets.Add($null, $n4)
$n6 = $wb.Worksheets.Add($null, $n5)

# --- Step 2: reposition so that the original sheet ends up 3rd (it will become "Sheet3") ---
# (worksheet object references above are position-bound, so re-fetch fresh by index for the moves)
$wb.Worksheets.Item(2).Move($wb.Worksheets.Item(1), $null)
$wb.Worksheets.Item(3).Move($wb.Worksheets.Item(2), $null)

# --- Step 3: rename sheets into final Sheet1..Sheet6 order ---
# (index3 currently still carries the name "Sheet1" from its original creation; free it up first)
$wb.Worksheets.Item(3).Name = "Sheet3_tmp"
$wb.Worksheets.Item(1).Name = "Sheet1"
$wb.Worksheets.Item(2).Name = "Sheet2"
$wb.Worksheets.Item(3).Name = "Sheet3"
$wb.Worksheets.Item(4).Name = "Sheet4"
$wb.Worksheets.Item(5).Name = "Sheet5"
$wb.Worksheets.Item(6).Name = "Sheet6"

# Sheet3 now holds the original content untouched; fill the rest with their target data.

# --- Sheet1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Value = "Search Data"
$ws.Range("B1").Value = "Titles To Verify"
$ws.Range("C1").Value = "Filter Verifier"
$ws.Range("D1").Value = "URL verifier"
$ws.Range("A2").Value = "Electronics"
$ws.Range("B2").Value = "Electronics"
$ws.Range("C2").Value = "Consumer Electronics"
$ws.Range("D2").Value = "/Manufacturing-Processing-Machinery-Catalog/Machine-Tools.html"
$ws.Range("C3").Value = "ISO 9000"
$ws.Range("D3").Value = "/Manufacturing-Processing-Machinery-Catalog/Engineering-Construction-Machinery.html"
$ws.Range("C4").Value = "Diamond Member"
$ws.Range("D4").Value = "/Manufacturing-Processing-Machinery-Catalog/Woodworking-Machinery.html"
$ws.Range("D5").Value = "/Manufacturing-Processing-Machinery-Catalog/Plastic-Machinery.html"
$ws.Range("D6").Value = "/Manufacturing-Processing-Machinery-Catalog/Metallic-Processing-Machinery.html"
$ws.Range("D7").Value = "/Manufacturing-Processing-Machinery-Catalog/Mould.html"
$ws.Range("D8").Value = "/Manufacturing-Processing-Machinery-Catalog/Laser-Equipment.html"
$ws.Range("D9").Value = "/Manufacturing-Processing-Machinery-Catalog/Casting-Forging.html"
$ws.Range("D10").Value = "/Manufacturing-Processing-Machinery-Catalog/Agricultural-Machinery.html"

# --- Sheet2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A1").Value = "Products"
$ws.Range("A2").Value = "Wires"
$ws.Range("A3").Value = "Furniture"

# --- Sheet4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("A1").Value = "Verify Title"
$ws.Range("B1").Value = "Verify label"
$ws.Range("A2").Value = "Register"
$ws.Range("B2").Value = "Español"
$ws.Range("A3").Value = "Sign In | Made-in-China.com"
$ws.Range("B3").Value = "Français"
$ws.Range("B4").Value = "Deutsch"
$ws.Range("B5").Value = "Русский язык"
$ws.Range("B6").Value = "日本語"
$ws.Range("B7").Value = "English"
$ws.Range("B8").Value = "Manufacturing & Processing Machinery"
$ws.Range("B9").Value = "Consumer Electronics"
$ws.Range("B10").Value = "Industrial Equipment & Components"
$ws.Range("B11").Value = "Electrical & Electronics"
$ws.Range("B12").Value = "Construction & Decoration"
$ws.Range("B13").Value = "Light Industry & Daily Use"
$ws.Range("B14").Value = "Auto, Motorcycle Parts & Accessories"
$ws.Range("B15").Value = "Apparel & Accessories"
$ws.Range("B16").Value = "Lights & Lighting"
$ws.Range("B17").Value = "Sporting Goods & Recreation"

# --- Sheet5 ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("A1").Value = "krushn34@gmail.com"
$ws.Range("A2").Value = "how-to-source-products-on-made-in-china-com"
$ws.Range("A3").Value = "audited-suppliers"
$ws.Range("A4").Value = "private-sourcing-meetings"

# --- Sheet6 ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("A1").Value = "LED lights"
$ws.Range("B1").Value = "Furniture.html"
$ws.Range("C1").Value = "Furniture"
$ws.Range("A2").Value = "I would like to inquire about your services. Please provide detailed information regarding pricing and availability at the earliest convenience"
$ws.Range("B2").Value = "Apparel-Clothing.html"
$ws.Range("C2").Value = "Apparel & Clothing"
$ws.Range("A3").Value = "tsr@gmail.com"
$ws.Range("B3").Value = "Auto-Parts-Accessories.html"
$ws.Range("C3").Value = "Auto Parts & Accessories"
$ws.Range("A4").Value = "john"
$ws.Range("B4").Value = "Bags-Cases-Luggages.html"
$ws.Range("C4").Value = "Bags, Cases & Luggages"
$ws.Range("A5").Value = "happy"
$ws.Range("B5").Value = "Computer-Products.html"
$ws.Range("C5").Value = "Computer Products"
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "1234567891"
$ws.Range("B6").Value = "Consumer-Electronics.html"
$ws.Range("C6").Value = "Consumer Electronics"
$ws.Range("B7").Value = "Electrical-Equipments.html"
$ws.Range("C7").Value = "Electrical Equipment"
$ws.Range("B8").Value = "Electronic-Components.html"
$ws.Range("C8").Value = "Electronic Components"
$ws.Range("B9").Value = "Gifts-Crafts-Collectibles.html"
$ws.Range("C9").Value = "Gifts, Crafts & Collectibles"
$ws.Range("B10").Value = "Motorcycles-Scooters.html"
$ws.Range("C10").Value = "Motorcycles & Scooters"
$ws.Range("B11").Value = "Sealing-Packaging-Storage-Shelving.html"
$ws.Range("C11").Value = "Sealing, Packaging, Storage & Shelving"
$ws.Range("B12").Value = "Building-Materials-Supplies.html"
$ws.Range("C12").Value = "Building Materials & Supplies"
